$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new regdNo values (column A) for rows 5-7
$ws.Range("A5").Value = "14221A0565"
$ws.Range("A6").Value = "14221A0568"
$ws.Range("A7").Value = "14221A0562"

# Fill in the new name values (column B) for rows 5-7
$ws.Range("B7").Value = "Aamir Shah"
$ws.Range("B6").Value = "Shailendra"
$ws.Range("B5").Value = "Pawan Aacharya"

# Fill in the new dues values (column C) for rows 5-7
$ws.Range("C5").Value = 78990
$ws.Range("C6").Value = 89000
$ws.Range("C7").Value = 9000

# Remove rows 8-13 (old leftover data)
$ws.Range("A8:C13").ClearContents()

# Set selection to C5 as the new active cell
$ws.Range("C5").Select()
